# Round anomalous times below 0.01s to 0.
# Applies to the "executionTime" worksheet: a specific set of small
# nonzero timing values (all previously rounded/displayed results of
# sub-0.01s measurements) are reset to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("executionTime")

$cells = @(
    "E12",
    "C16", "D16",
    "D17", "E17",
    "C18", "D18", "E18",
    "C20",
    "E22",
    "C23", "D23", "F23",
    "H26",
    "C28", "E28", "G28",
    "C30",
    "C34", "F34",
    "C38", "D38", "G38",
    "D43",
    "H47",
    "G48",
    "C49",
    "C51", "G51",
    "G52",
    "F53",
    "D54", "E54",
    "C57", "G57",
    "E58",
    "C59",
    "E60", "G60",
    "C61",
    "E63",
    "E68", "G68",
    "E69",
    "C73", "E73", "G73",
    "C76", "D76", "F76",
    "C87",
    "D90",
    "E91", "F91",
    "F96",
    "C99",
    "C101",
    "C102"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = 0
}
